# The workbook's first sheet ("sys_scr_Script") lists the mosaic scripts
# and their metadata, one row per script. Row 2 ("molgenis_mosaic") had its
# resultFileExtension (column E) set to "pdf"; update it to "mosaic.pdf".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sys_scr_Script")

$ws.Range("E2").Value = "mosaic.pdf"

# Leave the selection on the edited cell, matching the saved view state.
$ws.Range("E2").Select()
